# "first push of framework" - populate the Odoo test-data workbook:
#  - Sheet1 header row gets the full data1..data10 column set, styled with
#    a centred, bordered yellow header.
#  - Sheet1 rows 2-4 get the create/delete-customer test rows (values +
#    mailto hyperlinks on the email / address columns).
#  - Sheet2 gets a small "CustomersInfo" header + two rows.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Sheet1
# ---------------------------------------------------------------------

# Header row
$ws1.Range("A1:K1").Value = @("testcaseID","data1","data2","data3","data4","data5","data6","data7","data8","data9","data10")

# Row 2 - existing validLogin_ID row, new email + phone number
$ws1.Range("A2").Value = "validLogin_ID"
$ws1.Range("B2").Value = "pratap.ganesh27@gmail.com"
$ws1.Range("C2").Value = 9755314363

# Row 3 - create_customer_ID
$ws1.Range("A3").Value = "create_customer_ID"
$ws1.Range("B3").Value = "pratap.ganesh27@gmail.com"
$ws1.Range("C3").Value = 9755314363
$ws1.Range("D3").Value = "Ganesh"
$ws1.Range("E3").Value = "Sector-14"
$ws1.Range("F3").Value = "NehruPlace"
$ws1.Range("G3").Value = "Delhi"
$ws1.Range("H3").Value = "India"
$ws1.Range("I3").Value = 201265
$ws1.Range("J3").Value = 2415789944
$ws1.Range("K3").Value = "Meera.Chopra@gamil.com"

# Row 4 - delete_customer_ID
$ws1.Range("A4").Value = "delete_customer_ID"
$ws1.Range("B4").Value = "pratap.ganesh27@gmail.com"
$ws1.Range("C4").Value = 9755314363
$ws1.Range("D4").Value = "Meera"

# Header formatting - yellow fill, border, centred text
$header = $ws1.Range("A1:K1")
$header.Interior.Color = 65535
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4108

# Data-row borders (rows 2-18 already carry the bordered style from the
# template; re-assert it on the newly written rows for safety)
$ws1.Range("A2:K4").Borders.LineStyle = 1

# Hyperlinks - drop the old one and add the new set in the order the
# sheet references them (K3, B2, B3, B4)
$ws1.Cells.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("K3"), "mailto:Meera.Chopra@gamil.com")
$ws1.Hyperlinks.Add($ws1.Range("B2"), "mailto:pratap.ganesh27@gmail.com")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "mailto:pratap.ganesh27@gmail.com")
$ws1.Hyperlinks.Add($ws1.Range("B4"), "mailto:pratap.ganesh27@gmail.com")

# Column widths (characters, best-effort AutoFit equivalents)
$ws1.Columns.Item(1).ColumnWidth = 17.05
$ws1.Columns.Item(2).ColumnWidth = 21.83
$ws1.Columns.Item(3).ColumnWidth = 10.17
$ws1.Columns.Item(6).ColumnWidth = 9.5
$ws1.Columns.Item(10).ColumnWidth = 10.17
$ws1.Columns.Item(11).ColumnWidth = 22.5

# ---------------------------------------------------------------------
# Sheet2 - small CustomersInfo reference sheet
# ---------------------------------------------------------------------

$ws2.Range("A1").Value = "CustomersInfo"
$ws2.Range("B1").Value = "data1"
$ws2.Range("C1").Value = "data2"
$ws2.Range("D1").Value = "data3"
$ws2.Range("A1:I1").Interior.Color = 65535

$ws2.Range("A2").Value = "name"
$ws2.Range("A3").Value = "address"

$ws2.Columns.Item(1).ColumnWidth = 12.27

# Leave the selection on Sheet2 where the author left it, then reactivate
# Sheet1 (the originally tab-selected sheet) with its own selection.
$ws2.Range("C3").Select() | Out-Null
$ws1.Activate()
$ws1.Range("E4").Select() | Out-Null
